$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the "NB" row (old row 8); "SVM" (old row 9) shifts up to row 8.
$ws.Rows(8).Delete()

# 2) Insert 5 new blank columns, one right after each existing
#    "<period> Alt" column, to hold the new "<period> Alt std" data.
#    Working right-to-left keeps each letter valid at the moment of use.
$ws.Columns("H").Insert()
$ws.Columns("G").Insert()
$ws.Columns("F").Insert()
$ws.Columns("E").Insert()
$ws.Columns("D").Insert()

# 3) Relabel the header row: each old "<period> Alt" column becomes
#    "<period> Alt mean" and the newly inserted column next to it becomes
#    "<period> Alt std".
$ws.Range("C1").Value = "One Year Alt mean"
$ws.Range("D1").Value = "One Year Alt std"
$ws.Range("E1").Value = "Two Year Alt mean"
$ws.Range("F1").Value = "Two Year Alt std"
$ws.Range("G1").Value = "Three Year Alt mean"
$ws.Range("H1").Value = "Three Year Alt std"
$ws.Range("I1").Value = "Five Year Alt mean"
$ws.Range("J1").Value = "Five Year Alt std"
$ws.Range("K1").Value = "Ten Year Alt mean"
$ws.Range("L1").Value = "Ten Year Alt std"

# 4) Rename the "CART" algorithm row (row 5) to "DTREE".
$ws.Range("B5").Value = "DTREE"

# 4b) Row 8 ("SVM", shifted up from the old row 9) kept its old A-column
#     index (7) after the row delete above; fix it back to 6 so the
#     0-based index column stays contiguous (0..6).
$ws.Range("A8").Value = 6

# 5) Overwrite the mean columns and fill the new std columns with the
#    updated numbers for every remaining algorithm row (LR, LDA, KNN,
#    DTREE, RTREE, XTREE, SVM).
$ws.Range("C2").Value = 0.9059736235798841
$ws.Range("D2").Value = 0.00711537158413927
$ws.Range("E2").Value = 0.8996791978242232
$ws.Range("F2").Value = 0.01164431962057356
$ws.Range("G2").Value = 0.8972876073138771
$ws.Range("H2").Value = 0.01586012749448935
$ws.Range("I2").Value = 0.8924851062610987
$ws.Range("J2").Value = 0.0154495749963656
$ws.Range("K2").Value = 0.8939424273272254
$ws.Range("L2").Value = 0.007957667679034016
$ws.Range("C3").Value = 0.9072747014115092
$ws.Range("D3").Value = 0.00905540953701427
$ws.Range("E3").Value = 0.8986685290226607
$ws.Range("F3").Value = 0.01242082520910338
$ws.Range("G3").Value = 0.8937917161647458
$ws.Range("H3").Value = 0.01388236557267272
$ws.Range("I3").Value = 0.8875680242882511
$ws.Range("J3").Value = 0.01674382767960443
$ws.Range("K3").Value = 0.8925138558986541
$ws.Range("L3").Value = 0.01128968171563637
$ws.Range("C4").Value = 0.8709952066947381
$ws.Range("D4").Value = 0.01004177487076623
$ws.Range("E4").Value = 0.8737149459746425
$ws.Range("F4").Value = 0.01139211010137904
$ws.Range("G4").Value = 0.8761123290019963
$ws.Range("H4").Value = 0.01580657835323751
$ws.Range("I4").Value = 0.8750726785816578
$ws.Range("J4").Value = 0.01671271091790429
$ws.Range("K4").Value = 0.8860886777513857
$ws.Range("L4").Value = 0.009772689361251341
$ws.Range("C5").Value = 0.8890519318873971
$ws.Range("D5").Value = 0.007185834774605518
$ws.Range("E5").Value = 0.8784355641859858
$ws.Range("F5").Value = 0.01219181576947925
$ws.Range("G5").Value = 0.8804936744516431
$ws.Range("H5").Value = 0.009915466456784821
$ws.Range("I5").Value = 0.8809374462966145
$ws.Range("J5").Value = 0.0112322049635825
$ws.Range("K5").Value = 0.8939350752177354
$ws.Range("L5").Value = 0.02073730548299743
$ws.Range("C6").Value = 0.8949114165408755
$ws.Range("D6").Value = 0.00860484306114285
$ws.Range("E6").Value = 0.8804608763293418
$ws.Range("F6").Value = 0.01086036593397243
$ws.Range("G6").Value = 0.8643898570781232
$ws.Range("H6").Value = 0.01360163480075617
$ws.Range("I6").Value = 0.8582241364495617
$ws.Range("J6").Value = 0.01567532135476022
$ws.Range("K6").Value = 0.8470953512046149
$ws.Range("L6").Value = 0.01740784161532815
$ws.Range("C7").Value = 0.919150181404094
$ws.Range("D7").Value = 0.007877274826622909
$ws.Range("E7").Value = 0.9072668790206733
$ws.Range("F7").Value = 0.01082425416504697
$ws.Range("G7").Value = 0.898335027494397
$ws.Range("H7").Value = 0.01011892128668629
$ws.Range("I7").Value = 0.8951341152546257
$ws.Range("J7").Value = 0.01564437495892368
$ws.Range("K7").Value = 0.8967961769030653
$ws.Range("L7").Value = 0.01430570472898387
$ws.Range("C8").Value = 0.9079259023860597
$ws.Range("D8").Value = 0.008730573768285992
$ws.Range("E8").Value = 0.9030532985844959
$ws.Range("F8").Value = 0.00852133568105664
$ws.Range("G8").Value = 0.9009644471115574
$ws.Range("H8").Value = 0.01485101149404087
$ws.Range("I8").Value = 0.8981647763074985
$ws.Range("J8").Value = 0.0140319865552274
$ws.Range("K8").Value = 0.8963211175206423
$ws.Range("L8").Value = 0.01024910528281289
